# Auto-generated edit script: updates numeric cell values across all 8 sheets
# per the authoritative diff (commit: 'chore: update Sheets via scheduled runner').
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 102.15625
$ws.Range("I11").Value = 102.15625
$ws.Range("K11").Value = 102.15625
$ws.Range("M11").Value = 37.84375
$ws.Range("H98").Value = 3231.1
$ws.Range("J98").Value = 4160.875
$ws.Range("L98").Value = 4160.875
$ws.Range("N98").Value = -7156.875
$ws.Range("H111").Value = 1189.1428
$ws.Range("I111").Value = 986.875
$ws.Range("J111").Value = 1458.8334
$ws.Range("K111").Value = 2960.625
$ws.Range("L111").Value = 4376.5002
$ws.Range("M111").Value = 106.375
$ws.Range("N111").Value = -10510.5002
$ws.Range("H116").Value = 6226.294
$ws.Range("J116").Value = 5763.7144
$ws.Range("L116").Value = 5763.7144
$ws.Range("N116").Value = -12647.7144
$ws.Range("H122").Value = 3231.1
$ws.Range("J122").Value = 4160.875
$ws.Range("L122").Value = 12482.625
$ws.Range("N122").Value = -17382.625
$ws.Range("H137").Value = 16246.857
$ws.Range("I137").Value = 1404.5454
$ws.Range("K137").Value = 4213.6362
$ws.Range("M137").Value = -1663.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 3041.1667
$ws.Range("J2").Value = 1622.5
$ws.Range("K2").Value = 3041.1667
$ws.Range("L2").Value = 1622.5
$ws.Range("M2").Value = -2928.1667
$ws.Range("N2").Value = -1848.5
$ws.Range("H32").Value = 3178670.8
$ws.Range("I32").Value = 4204502
$ws.Range("K32").Value = 4204502
$ws.Range("M32").Value = -4204215
$ws.Range("H45").Value = 1528.64
$ws.Range("I45").Value = 1500.8
$ws.Range("K45").Value = 1500.8
$ws.Range("M45").Value = -1123.8
$ws.Range("H61").Value = 91006.586
$ws.Range("I61").Value = 3307.2
$ws.Range("K61").Value = 3307.2
$ws.Range("M61").Value = -3095.2
$ws.Range("I116").Value = 3041.1667
$ws.Range("J116").Value = 1622.5
$ws.Range("K116").Value = 3041.1667
$ws.Range("L116").Value = 1622.5
$ws.Range("M116").Value = -747.1667000000002
$ws.Range("N116").Value = -6210.5
$ws.Range("H122").Value = 4148183
$ws.Range("I122").Value = 4148183
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12444549
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12442099
$ws.Range("N122").ClearContents()
$ws.Range("H125").Value = 92000
$ws.Range("J125").Value = 92000
$ws.Range("L125").Value = 92000
$ws.Range("N125").Value = -101840
$ws.Range("H132").Value = 2635709.2
$ws.Range("I132").Value = 1429.9656
$ws.Range("K132").Value = 4289.8968
$ws.Range("M132").Value = -1759.8968
$ws.Range("H136").Value = 91006.586
$ws.Range("I136").Value = 3307.2
$ws.Range("K136").Value = 9921.599999999999
$ws.Range("M136").Value = -7371.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 3041.1667
$ws.Range("J3").Value = 1622.5
$ws.Range("K3").Value = 3041.1667
$ws.Range("L3").Value = 1622.5
$ws.Range("M3").Value = -2927.1667
$ws.Range("N3").Value = -1850.5
$ws.Range("H105").Value = 43479564
$ws.Range("I105").Value = 47620364
$ws.Range("K105").Value = 47620364
$ws.Range("M105").Value = -47618617
$ws.Range("H134").Value = 29129.125
$ws.Range("I134").Value = 30519.473
$ws.Range("J134").Value = 24958.084
$ws.Range("K134").Value = 91558.41900000001
$ws.Range("L134").Value = 74874.25199999999
$ws.Range("M134").Value = -89023.41900000001
$ws.Range("N134").Value = -79944.25199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 33335478
$ws.Range("I132").Value = 2253.6667
$ws.Range("J132").Value = 333334500
$ws.Range("K132").Value = 6761.000100000001
$ws.Range("L132").Value = 1000003500
$ws.Range("M132").Value = -4231.000100000001
$ws.Range("N132").Value = -1000008560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 3499
$ws.Range("I87").Value = 3499
$ws.Range("K87").Value = 10497
$ws.Range("M87").Value = -9249
$ws.Range("H90").Value = 3499
$ws.Range("I90").Value = 3499
$ws.Range("K90").Value = 31491
$ws.Range("M90").Value = -25251
$ws.Range("H138").Value = 3135.7778
$ws.Range("I138").Value = 2746
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 8238
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = -3098
$ws.Range("N138").Value = -23780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2295.875
$ws.Range("I43").Value = 2295.875
$ws.Range("K43").Value = 2295.875
$ws.Range("M43").Value = -2144.875
$ws.Range("H103").Value = 102500
$ws.Range("J103").Value = 102500
$ws.Range("L103").Value = 102500
$ws.Range("N103").Value = -104844
$ws.Range("H113").Value = 11999.444
$ws.Range("I113").Value = 5997
$ws.Range("J113").Value = 16801.4
$ws.Range("K113").Value = 5997
$ws.Range("L113").Value = 16801.4
$ws.Range("M113").Value = -3827
$ws.Range("N113").Value = -21141.4
$ws.Range("H122").Value = 2832999.8
$ws.Range("I122").Value = 3089949.5
$ws.Range("J122").Value = 6553
$ws.Range("K122").Value = 9269848.5
$ws.Range("L122").Value = 19659
$ws.Range("M122").Value = -9267398.5
$ws.Range("N122").Value = -24559
$ws.Range("H126").Value = 8261624
$ws.Range("J126").Value = 11908669
$ws.Range("L126").Value = 35726007
$ws.Range("N126").Value = -35730947
$ws.Range("H132").Value = 9334.392
$ws.Range("I132").Value = 5458.4116
$ws.Range("J132").Value = 20316.334
$ws.Range("K132").Value = 16375.2348
$ws.Range("L132").Value = 60949.00199999999
$ws.Range("M132").Value = -13845.2348
$ws.Range("N132").Value = -66009.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5967824.5
$ws.Range("I7").Value = 9545410
$ws.Range("K7").Value = 9545410
$ws.Range("M7").Value = -9545298
$ws.Range("H16").Value = 125001840
$ws.Range("I16").Value = 142859020
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 142859020
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -142858850
$ws.Range("N16").Value = -1840
$ws.Range("H22").Value = 125001530
$ws.Range("J22").Value = 333334080
$ws.Range("L22").Value = 333334080
$ws.Range("N22").Value = -333334670
$ws.Range("H27").Value = 125001530
$ws.Range("J27").Value = 333334080
$ws.Range("L27").Value = 333334080
$ws.Range("N27").Value = -333334294
$ws.Range("H55").Value = 1633.9131
$ws.Range("I55").Value = 1801.2222
$ws.Range("J55").Value = 1526.3572
$ws.Range("K55").Value = 1801.2222
$ws.Range("L55").Value = 1526.3572
$ws.Range("M55").Value = -1628.2222
$ws.Range("N55").Value = -1872.3572
$ws.Range("H68").Value = 11543
$ws.Range("I68").Value = 18496.5
$ws.Range("J68").Value = 3198.8
$ws.Range("K68").Value = 18496.5
$ws.Range("L68").Value = 3198.8
$ws.Range("M68").Value = -17747.5
$ws.Range("N68").Value = -4696.8
$ws.Range("H71").Value = 11543
$ws.Range("I71").Value = 18496.5
$ws.Range("J71").Value = 3198.8
$ws.Range("K71").Value = 92482.5
$ws.Range("L71").Value = 15994
$ws.Range("M71").Value = -88738.5
$ws.Range("N71").Value = -23482
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H126").Value = 5967824.5
$ws.Range("I126").Value = 9545410
$ws.Range("K126").Value = 28636230
$ws.Range("M126").Value = -28633760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 20000
$ws.Range("J97").Value = 20000
$ws.Range("L97").Value = 20000
$ws.Range("N97").Value = -21982
$ws.Range("H113").Value = 5620.1665
$ws.Range("I113").Value = 5444.2
$ws.Range("K113").Value = 16332.6
$ws.Range("M113").Value = -14162.6
$ws.Range("H132").Value = 16604.625
$ws.Range("I132").Value = 3023.5
$ws.Range("J132").Value = 30185.75
$ws.Range("K132").Value = 9070.5
$ws.Range("L132").Value = 90557.25
$ws.Range("M132").Value = -6540.5
$ws.Range("N132").Value = -95617.25
